$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "47÷9=5, 2"
$t.Cell(1,2).Range.Text = "12÷8=1, 4"
$t.Cell(1,3).Range.Text = "83÷5=16, 3"
$t.Cell(1,4).Range.Text = "30÷7=4, 2"
$t.Cell(1,5).Range.Text = "55÷6=9, 1"
$t.Cell(5,1).Range.Text = "71÷5=14, 1"
$t.Cell(5,2).Range.Text = "81÷4=20, 1"
$t.Cell(5,3).Range.Text = "17÷8=2, 1"
$t.Cell(5,4).Range.Text = "40÷9=4, 4"
$t.Cell(5,5).Range.Text = "72÷4=18, 0"
$t.Cell(9,1).Range.Text = "28÷7=4, 0"
$t.Cell(9,2).Range.Text = "87÷9=9, 6"
$t.Cell(9,3).Range.Text = "61÷8=7, 5"
$t.Cell(9,4).Range.Text = "59÷8=7, 3"
$t.Cell(9,5).Range.Text = "14÷2=7, 0"
$t.Cell(13,1).Range.Text = "45÷3=15, 0"
$t.Cell(13,2).Range.Text = "31÷8=3, 7"
$t.Cell(13,3).Range.Text = "13÷2=6, 1"
$t.Cell(13,4).Range.Text = "47÷4=11, 3"
$t.Cell(13,5).Range.Text = "16÷8=2, 0"
$t.Cell(17,1).Range.Text = "21÷2=10, 1"
$t.Cell(17,2).Range.Text = "28÷8=3, 4"
$t.Cell(17,3).Range.Text = "47÷4=11, 3"
$t.Cell(17,4).Range.Text = "14÷5=2, 4"
$t.Cell(17,5).Range.Text = "97÷5=19, 2"
